$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the existing
# header formatting (bold, bordered, centered) from B1:C1 so the new
# headers share the same style as the rest of the header row.
$ws.Range("B1:C1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J numeric columns for rows 2-15.
$data = @(
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(4, 5),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(7, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
